$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers for new columns I and J ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy formatting (bold font, border, center/top alignment) from an
# existing header cell (H1) onto the two new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Data values for columns I ("I0") and J ("IF"), rows 2-31 ---
$values = @{
    2  = @(1, 6)
    3  = @(1, 7)
    4  = @(1, 6)
    5  = @(1, 7)
    6  = @(1, 6)
    7  = @(1, 5)
    8  = @(1, 7)
    9  = @(1, 6)
    10 = @(1, 6)
    11 = @(1, 8)
    12 = @(1, 6)
    13 = @(1, 7)
    14 = @(1, 7)
    15 = @(1, 6)
    16 = @(1, 7)
    17 = @(1, 5)
    18 = @(1, 9)
    19 = @(1, 6)
    20 = @(1, 5)
    21 = @(1, 5)
    22 = @(1, 5)
    23 = @(1, 6)
    24 = @(1, 7)
    25 = @(1, 7)
    26 = @(8, 9)
    27 = @(7, 7)
    28 = @(7, 9)
    29 = @(5, 7)
    30 = @(1, 3)
    31 = @(1, 3)
}

foreach ($row in 2..31) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
